$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column widths for the new columns ----
$ws.Columns.Item(3).ColumnWidth = 31.140625
$ws.Columns.Item(4).ColumnWidth = 16
$ws.Columns.Item(5).ColumnWidth = 15.140625
$ws.Columns.Item(6).ColumnWidth = 13.85546875
$ws.Columns.Item(7).ColumnWidth = 16.85546875
$ws.Columns.Item(8).ColumnWidth = 11.28515625

# ---- Header row (row 1) ----
$ws.Range("C1").Value = "store_name"
$ws.Range("D1").Value = "brand_name"
$ws.Range("E1").Value = "corporate_id"
$ws.Range("F1").Value = "store_address"
$ws.Range("G1").Value = "phoneNumber"
$ws.Range("H1").Value = "store_count"

# ---- Row 2 : Cricket Wireless ----
$ws.Range("C2").Value = "Cricket Wireless Authorized Retailer"
$ws.Range("D2").Value = "Cricket Wireless"
$ws.Range("E2").Value = 70189192
$ws.Range("F2").Value = "365 N Saddle Creek Rd"
$ws.Range("G2").Value = "(402) 885-6815"
$ws.Range("H2").Value = 50

# ---- Row 3 : Dippin' Dots ----
$ws.Range("C3").Value = "Dippin' Dots"
$ws.Range("D3").Value = "Dippin' Dots"
$ws.Range("E3").Value = 61
$ws.Range("F3").Value = "10000 California St."
$ws.Range("G3").Value = "402-393-0663"
$ws.Range("H3").Value = 30

# ---- Row 4 : Dollar Rent A Car ----
$ws.Range("C4").Value = "Dollar Rent A Car"
$ws.Range("D4").Value = "Dollar Rent A Car"
$ws.Range("E4").Value = 9346132
$ws.Range("F4").Value = "4501 ABBOT DRIVE"
$ws.Range("G4").Value = "(402) 345-2783"
$ws.Range("H4").Value = 28

# ---- Row 5 : WaBa Grill ----
$ws.Range("C5").Value = "WaBa Grill"
$ws.Range("D5").Value = "WaBa Grill"
$ws.Range("E5").Value = 9002
$ws.Range("F5").Value = "13131 Crossroads Pkway S."
$ws.Range("G5").Value = "562-463-9222"
$ws.Range("H5").Value = 29

# ---- Row 6 : Arvest Bank ----
$ws.Range("C6").Value = "Arvest Bank"
$ws.Range("D6").Value = "Arvest Bank"
$ws.Range("E6").Value = 964
$ws.Range("F6").Value = "9221 North Oak Trafficway"
$ws.Range("G6").Value = "(913) 279-3300"
$ws.Range("H6").Value = 12

# ---- Borders: box (all 4 sides, thin) around every new cell C1:H6 ----
$box = $ws.Range("C1:H6")
$box.Borders.LineStyle = 1
$box.Borders.Weight = 2

# ---- store_count header (H1) only needs a left/right border, no top/bottom ----
$ws.Range("H1").Borders(8).LineStyle = -4142
$ws.Range("H1").Borders(9).LineStyle = -4142

# ---- F4 (address for Dollar Rent A Car) also only has a left/right border ----
$ws.Range("F4").Borders(8).LineStyle = -4142
$ws.Range("F4").Borders(9).LineStyle = -4142

# ---- Header row fill: match the yellow header look for store_count ----
$ws.Range("H1").Interior.Color = 65535

# ---- Number format: corporate_id column stored as text ----
$ws.Range("E1:E6").NumberFormat = "@"

# ---- Print setup ----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---- Selection ----
$ws.Range("G13").Select()
